$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append a new log row (row 67) with the latest run information.
$newRow = 67
$prevRow = $newRow - 1

# Carry the formatting of the previous log row down onto the new row.
$ws.Range("A$prevRow`:H$prevRow").Copy()
$ws.Range("A$newRow`:H$newRow").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item($newRow, 1).Value = "2025-08-28 06:48:04 UTC"
$ws.Cells.Item($newRow, 2).Value = "2025-08-28 12:18:04 IST"
$ws.Cells.Item($newRow, 3).Value = "UPDATED"
$ws.Cells.Item($newRow, 4).Value = "New circular processed."
$ws.Cells.Item($newRow, 5).Value = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-28-08-2025.pdf"
$ws.Cells.Item($newRow, 6).Value = "INGOT-28-08-2025.pdf"
$ws.Cells.Item($newRow, 7).Value = 1
$ws.Cells.Item($newRow, 8).Value = 4
